$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("List1")
$ws2 = $wb.Worksheets.Item("Měření aktivity")

# --- Sheet1 (List1) changes ---
$ws1.Activate()

# Update A19 value (date/time) - this ripples into B19/C19/E18/C20 formulas
$ws1.Range("A19").Value = 45498.78402777778

# Cursor/selection moved to A20
$ws1.Range("A20").Select()

# --- Sheet2 (Měření aktivity) changes ---
$ws2.Activate()

# Corrected A17 value
$ws2.Range("A17").Value = 45496.621527777781

# Fill in the previously-blank measurement row 18
$ws2.Range("A18").Value = 45498.78402777778
$ws2.Range("B18").Value = 0.004
$ws2.Range("C18").Value = 0.004
$ws2.Range("D18").Value = 0.004
$ws2.Range("E18").Value = 0.003
$ws2.Range("F18").Value = 0.005
$ws2.Range("G18").Value = 158.7
$ws2.Range("H18").Value = 158.6
$ws2.Range("I18").Value = 158.5
$ws2.Range("J18").Value = 158.5
$ws2.Range("K18").Value = 158.5
$ws2.Range("L18").Value = 158.6
$ws2.Range("M18").Value = 158.5
$ws2.Range("N18").Value = 158.5
$ws2.Range("O18").Value = 158.5
$ws2.Range("P18").Value = 158.5
$ws2.Range("Q18").Formula = "=AVERAGE(G18:P18) - AVERAGE(Tabulka1[[#This Row],[č.1]:[č.5]])"

# Cursor/selection moved to Q18
$ws2.Range("Q18").Select()

$excel.CalculateFull()

